$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 19.69872420814777
$ws.Cells.Item(2, 3).Value = 10.02792530156427
$ws.Cells.Item(2, 4).Value = 8.762349761672796
$ws.Cells.Item(2, 6).Value = 35.17601278616402
$ws.Cells.Item(2, 7).Value = 3.671341164948863
$ws.Cells.Item(2, 10).Value = 10.4517895089887
$ws.Cells.Item(2, 12).Value = 11.79906010375989
$ws.Cells.Item(2, 15).Value = 26.7331436993317

# Row 3
$ws.Cells.Item(3, 2).Value = 19.14949602794881
$ws.Cells.Item(3, 3).Value = 9.776263997931203
$ws.Cells.Item(3, 4).Value = 8.758352327594721
$ws.Cells.Item(3, 6).Value = 35.31220929543376
$ws.Cells.Item(3, 7).Value = 3.673688856092156
$ws.Cells.Item(3, 10).Value = 10.48471291032264
$ws.Cells.Item(3, 12).Value = 11.77109248581744
$ws.Cells.Item(3, 15).Value = 26.85708774456107

# Row 4
$ws.Cells.Item(4, 2).Value = 18.80584932469516
$ws.Cells.Item(4, 3).Value = 9.617961209899061
$ws.Cells.Item(4, 4).Value = 8.756786406887691
$ws.Cells.Item(4, 6).Value = 35.40552805979562
$ws.Cells.Item(4, 7).Value = 3.67520657776264
$ws.Cells.Item(4, 10).Value = 10.50598998336855
$ws.Cells.Item(4, 12).Value = 11.75524228523777
$ws.Cells.Item(4, 15).Value = 26.9402314876485

# Row 5
$ws.Cells.Item(5, 2).Value = 18.66439295784122
$ws.Cells.Item(5, 3).Value = 9.552572284319423
$ws.Cells.Item(5, 4).Value = 8.756372743332953
$ws.Cells.Item(5, 6).Value = 35.44598479904298
$ws.Cells.Item(5, 7).Value = 3.675844292555838
$ws.Cells.Item(5, 10).Value = 10.51492836855662
$ws.Cells.Item(5, 12).Value = 11.74911950815474
$ws.Cells.Item(5, 15).Value = 26.97587826497966

# Row 6
$ws.Cells.Item(6, 2).Value = 18.64082460688844
$ws.Cells.Item(6, 3).Value = 9.541663581706779
$ws.Cells.Item(6, 4).Value = 8.75631764081556
$ws.Cells.Item(6, 6).Value = 35.45284902337908
$ws.Cells.Item(6, 7).Value = 3.675951347941727
$ws.Cells.Item(6, 10).Value = 10.51642877785702
$ws.Cells.Item(6, 12).Value = 11.74812324663826
$ws.Cells.Item(6, 15).Value = 26.98190382270185

# Row 7
$ws.Cells.Item(7, 2).Value = 18.80394706544913
$ws.Cells.Item(7, 3).Value = 9.617082814003236
$ws.Cells.Item(7, 4).Value = 8.756779917907501
$ws.Cells.Item(7, 6).Value = 35.40606385235465
$ws.Cells.Item(7, 7).Value = 3.675215100259798
$ws.Cells.Item(7, 10).Value = 10.50610944421201
$ws.Cells.Item(7, 12).Value = 11.75515834467256
$ws.Cells.Item(7, 15).Value = 26.94070509355576

# Row 8
$ws.Cells.Item(8, 2).Value = 19.51079633410599
$ws.Cells.Item(8, 3).Value = 9.941978981616536
$ws.Cells.Item(8, 4).Value = 8.760787672126121
$ws.Cells.Item(8, 6).Value = 35.22095615842319
$ws.Cells.Item(8, 7).Value = 3.672134864916737
$ws.Cells.Item(8, 10).Value = 10.462921546515
$ws.Cells.Item(8, 12).Value = 11.78914486831413
$ws.Cells.Item(8, 15).Value = 26.77441474858544

# Row 9
$ws.Cells.Item(9, 2).Value = 20.83803209542954
$ws.Cells.Item(9, 3).Value = 10.54617866943292
$ws.Cells.Item(9, 4).Value = 8.775649834543227
$ws.Cells.Item(9, 6).Value = 34.9352618417667
$ws.Cells.Item(9, 7).Value = 3.666696538448063
$ws.Cells.Item(9, 10).Value = 10.38662144906004
$ws.Cells.Item(9, 12).Value = 11.86610396282755
$ws.Cells.Item(9, 15).Value = 26.5044417120294

# Row 10
$ws.Cells.Item(10, 2).Value = 21.76753622624952
$ws.Cells.Item(10, 3).Value = 10.96644761282133
$ws.Cells.Item(10, 4).Value = 8.790772907282072
$ws.Cells.Item(10, 6).Value = 34.7729977176269
$ws.Cells.Item(10, 7).Value = 3.663064015715037
$ws.Cells.Item(10, 10).Value = 10.3356300131896
$ws.Cells.Item(10, 12).Value = 11.92868209981742
$ws.Cells.Item(10, 15).Value = 26.34064426097234

# Row 11
$ws.Cells.Item(11, 2).Value = 22.17870610781431
$ws.Cells.Item(11, 3).Value = 11.15184690790492
$ws.Cells.Item(11, 4).Value = 8.798550413611387
$ws.Cells.Item(11, 6).Value = 34.70962522628148
$ws.Cells.Item(11, 7).Value = 3.661489464422273
$ws.Cells.Item(11, 10).Value = 10.31352229283561
$ws.Cells.Item(11, 12).Value = 11.95840372920931
$ws.Cells.Item(11, 15).Value = 26.27370247188195

# Row 12
$ws.Cells.Item(12, 2).Value = 22.33259741424586
$ws.Cells.Item(12, 3).Value = 11.22117235397011
$ws.Cells.Item(12, 4).Value = 8.801623251632329
$ws.Cells.Item(12, 6).Value = 34.68713707395631
$ws.Cells.Item(12, 7).Value = 3.660904360570963
$ws.Cells.Item(12, 10).Value = 10.30530641934412
$ws.Cells.Item(12, 12).Value = 11.96983370434242
$ws.Cells.Item(12, 15).Value = 26.24944781277919

# Row 13
$ws.Cells.Item(13, 2).Value = 22.29953662827278
$ws.Cells.Item(13, 3).Value = 11.20628181485737
$ws.Cells.Item(13, 4).Value = 8.800955806979623
$ws.Cells.Item(13, 6).Value = 34.69191303515021
$ws.Cells.Item(13, 7).Value = 3.661029878410009
$ws.Cells.Item(13, 10).Value = 10.30706893458167
$ws.Cells.Item(13, 12).Value = 11.96736435711193
$ws.Cells.Item(13, 15).Value = 26.25462270221808

# Row 14
$ws.Cells.Item(14, 2).Value = 22.19140382103389
$ws.Cells.Item(14, 3).Value = 11.15756828381844
$ws.Cells.Item(14, 4).Value = 8.79880066783547
$ws.Cells.Item(14, 6).Value = 34.70774481262784
$ws.Cells.Item(14, 7).Value = 3.661441104551055
$ws.Cells.Item(14, 10).Value = 10.31284324863654
$ws.Cells.Item(14, 12).Value = 11.95934060346414
$ws.Cells.Item(14, 15).Value = 26.27168505150358

# Row 15
$ws.Cells.Item(15, 2).Value = 22.12492996010111
$ws.Cells.Item(15, 3).Value = 11.12761367939701
$ws.Cells.Item(15, 4).Value = 8.79749716422914
$ws.Cells.Item(15, 6).Value = 34.71763906367185
$ws.Cells.Item(15, 7).Value = 3.661694442084249
$ws.Cells.Item(15, 10).Value = 10.3164004573666
$ws.Cells.Item(15, 12).Value = 11.9544484556549
$ws.Cells.Item(15, 15).Value = 26.28227898089546

# Row 16
$ws.Cells.Item(16, 2).Value = 21.7404189557145
$ws.Cells.Item(16, 3).Value = 10.95421061500774
$ws.Cells.Item(16, 4).Value = 8.790282574791162
$ws.Cells.Item(16, 6).Value = 34.77735008240088
$ws.Cells.Item(16, 7).Value = 3.663168478989808
$ws.Cells.Item(16, 10).Value = 10.33709664975757
$ws.Cells.Item(16, 12).Value = 11.92676451181705
$ws.Cells.Item(16, 15).Value = 26.34517192696895

# Row 17
$ws.Cells.Item(17, 2).Value = 21.50144771260146
$ws.Cells.Item(17, 3).Value = 10.84631549170931
$ws.Cells.Item(17, 4).Value = 8.786085627028568
$ws.Cells.Item(17, 6).Value = 34.81666124221586
$ws.Cells.Item(17, 7).Value = 3.664092664226576
$ws.Cells.Item(17, 10).Value = 10.35007140142226
$ws.Cells.Item(17, 12).Value = 11.91009886956907
$ws.Cells.Item(17, 15).Value = 26.38569807848186

# Row 18
$ws.Cells.Item(18, 2).Value = 21.36290675829617
$ws.Cells.Item(18, 3).Value = 10.78371540669469
$ws.Cells.Item(18, 4).Value = 8.783756264063067
$ws.Cells.Item(18, 6).Value = 34.84025440360514
$ws.Cells.Item(18, 7).Value = 3.664631567005856
$ws.Cells.Item(18, 10).Value = 10.35763663476232
$ws.Cells.Item(18, 12).Value = 11.90063164753499
$ws.Cells.Item(18, 15).Value = 26.40971985935219

# Row 19
$ws.Cells.Item(19, 2).Value = 21.31581602984764
$ws.Cells.Item(19, 3).Value = 10.76242864878678
$ws.Cells.Item(19, 4).Value = 8.782982157014418
$ws.Cells.Item(19, 6).Value = 34.84841112213856
$ws.Cells.Item(19, 7).Value = 3.664815291933021
$ws.Cells.Item(19, 10).Value = 10.36021571828641
$ws.Cells.Item(19, 12).Value = 11.89744670680778
$ws.Cells.Item(19, 15).Value = 26.41797535658655

# Row 20
$ws.Cells.Item(20, 2).Value = 21.52700052414354
$ws.Cells.Item(20, 3).Value = 10.85785754579958
$ws.Cells.Item(20, 4).Value = 8.786523652782575
$ws.Cells.Item(20, 6).Value = 34.81237476538796
$ws.Cells.Item(20, 7).Value = 3.663993524313562
$ws.Cells.Item(20, 10).Value = 10.34867961367486
$ws.Cells.Item(20, 12).Value = 11.91186073839382
$ws.Cells.Item(20, 15).Value = 26.38131024472076

# Row 21
$ws.Cells.Item(21, 2).Value = 22.22321515525665
$ws.Cells.Item(21, 3).Value = 11.1719009186724
$ws.Cells.Item(21, 4).Value = 8.799430231603383
$ws.Cells.Item(21, 6).Value = 34.70305360284905
$ws.Cells.Item(21, 7).Value = 3.661320015462153
$ws.Cells.Item(21, 10).Value = 10.31114296863328
$ws.Cells.Item(21, 12).Value = 11.96169266754245
$ws.Cells.Item(21, 15).Value = 26.26664366663291

# Row 22
$ws.Cells.Item(22, 2).Value = 22.66763035041606
$ws.Cells.Item(22, 3).Value = 11.37198996593224
$ws.Cells.Item(22, 4).Value = 8.808608833478502
$ws.Cells.Item(22, 6).Value = 34.64040862068528
$ws.Cells.Item(22, 7).Value = 3.659637655418847
$ws.Cells.Item(22, 10).Value = 10.2875186424694
$ws.Cells.Item(22, 12).Value = 11.99527835962634
$ws.Cells.Item(22, 15).Value = 26.19808619931942

# Row 23
$ws.Cells.Item(23, 2).Value = 22.43144800554971
$ws.Cells.Item(23, 3).Value = 11.26568573574335
$ws.Cells.Item(23, 4).Value = 8.803642518255094
$ws.Cells.Item(23, 6).Value = 34.67303543608693
$ws.Cells.Item(23, 7).Value = 3.660529640226323
$ws.Cells.Item(23, 10).Value = 10.3000445305583
$ws.Cells.Item(23, 12).Value = 11.97726174761518
$ws.Cells.Item(23, 15).Value = 26.23409049956706

# Row 24
$ws.Cells.Item(24, 2).Value = 21.51545167983641
$ws.Cells.Item(24, 3).Value = 10.85264115398649
$ws.Cells.Item(24, 4).Value = 8.786325361094883
$ws.Cells.Item(24, 6).Value = 34.81430959006965
$ws.Cells.Item(24, 7).Value = 3.664038321863641
$ws.Cells.Item(24, 10).Value = 10.3493085110603
$ws.Cells.Item(24, 12).Value = 11.91106384176191
$ws.Cells.Item(24, 15).Value = 26.38329173342676

# Row 25
$ws.Cells.Item(25, 2).Value = 20.4863407426591
$ws.Cells.Item(25, 3).Value = 10.38663787905199
$ws.Cells.Item(25, 4).Value = 8.770885988177096
$ws.Cells.Item(25, 6).Value = 35.00421996796262
$ws.Cells.Item(25, 7).Value = 3.668103714263384
$ws.Cells.Item(25, 10).Value = 10.38662144906004
$ws.Cells.Item(25, 12).Value = 11.86610396282755
$ws.Cells.Item(25, 15).Value = 26.57143315949245
